# Generovani ok predmetu a zaklad pro xlsx -> csv
#
# Re-pairs several rows that share the same seminariciUcitIdno so that the
# "K"-prefixed (katedra) variant and the plain variant of zkratka swap
# places, and fixes the Pocitacove modelovani I / Programovani A ordering
# for idno 612.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Text {
    # Safe swap of two cells' text values (works for plain text strings).
    param($addr1, $addr2)
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

function Swap-TextSafe {
    # Swap two cells' text values without letting Excel coerce numeric-looking
    # strings (e.g. "0182") into actual numbers. Uses copy/paste-values so the
    # original text type and cell style are preserved exactly.
    param($addr1, $addr2)
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $tmpAddr = "ZZ1"
    $tmp = $ws.Range($tmpAddr)

    $r1.Copy() | Out-Null
    $tmp.PasteSpecial(-4163) | Out-Null   # xlPasteValues

    $r2.Copy() | Out-Null
    $r1.PasteSpecial(-4163) | Out-Null    # xlPasteValues

    $tmp.Copy() | Out-Null
    $r2.PasteSpecial(-4163) | Out-Null    # xlPasteValues

    $tmp.ClearContents() | Out-Null
}

# idno 14 (rows 2-3): swap EIS / KEIS
Swap-Text "B2" "B3"

# idno 612 (rows 8-10): Pocitacove modelovani I / Programovani A / Pocitacove modelovani I
# -> Pocitacove modelovani I / Pocitacove modelovani I / Programovani A
$ws.Range("A9").Value2 = "Počítačové modelování I"
$ws.Range("A10").Value2 = "Programování A"

$b8 = $ws.Range("B8").Value2
$b9 = $ws.Range("B9").Value2
$b10 = $ws.Range("B10").Value2
$ws.Range("B8").Value2 = $b10
$ws.Range("B9").Value2 = $b8
$ws.Range("B10").Value2 = $b9

# idno 1609 (rows 11-12): swap EIS / KEIS
Swap-Text "B11" "B12"

# idno 3457 (rows 14-15): swap EIS / KEIS
Swap-Text "B14" "B15"

# idno 3606 (rows 16-17): swap EIS / KEIS
Swap-Text "B16" "B17"

# idno 4746 (rows 21-22): swap Zaklady autonomni robotiky / Prakticke aplikace hardwaru
# and 0182 / AHW (0182 needs the numeric-safe swap so it stays text "0182")
Swap-Text "A21" "A22"
Swap-TextSafe "B21" "B22"

# idno 4991 (rows 24-25): swap EIS / KEIS
Swap-Text "B24" "B25"

# idno 8021 (rows 27-28): swap RSPP / KSPP
Swap-Text "B27" "B28"

# idno 8514 (rows 32, 34): swap OPRE / KOPRE
Swap-Text "B32" "B34"
